$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Dhh", "Hhip", "ECs", 3, 1, 0.9167996666666666, 2.750399, 0.2833456974325495, 0.2833456974325495, 3, 1, 0.3309496666666666, 0.9928489999999999, 0.02814339759365984, 0.02814339759365983, 0.3034145440834444, 2.730730896751, 0.007974310619297082, 0.007974310619297081),
    @("ECs", "Dhh", "Hhip", "FAPs", 3, 1, 0.9167996666666666, 2.750399, 0.2833456974325495, 0.2833456974325495, 3, 1, 6.605175333333333, 19.815526, 0.561692892620634, 0.561692892620634, 6.055622543874888, 54.50060289487399, 0.1591532644024997, 0.1591532644024997),
    @("ECs", "Dhh", "Hhip", "MuSCs", 3, 1, 0.9167996666666666, 2.750399, 0.2833456974325495, 0.2833456974325495, 3, 1, 4.823282, 14.469846, 0.410163709785706, 0.410163709785706, 4.421983329839333, 39.797849968554, 0.1162181224107527, 0.1162181224107527),
    @("FAPs", "Dhh", "Hhip", "ECs", 3, 1, 0.7348883333333333, 2.204665, 0.2271242616180895, 0.2271242616180895, 3, 1, 0.3309496666666666, 0.9928489999999999, 0.02814339759365984, 0.02814339759365983, 0.2432110489538888, 2.188899440585, 0.006392048397884307, 0.006392048397884307),
    @("FAPs", "Dhh", "Hhip", "FAPs", 3, 1, 0.7348883333333333, 2.204665, 0.2271242616180895, 0.2271242616180895, 3, 1, 6.605175333333333, 19.815526, 0.561692892620634, 0.561692892620634, 4.854066292087777, 43.68659662878999, 0.1275740834925903, 0.1275740834925904),
    @("FAPs", "Dhh", "Hhip", "MuSCs", 3, 1, 0.7348883333333333, 2.204665, 0.2271242616180895, 0.2271242616180895, 3, 1, 4.823282, 14.469846, 0.410163709785706, 0.410163709785706, 3.544573670176666, 31.90116303159, 0.09315812972761484, 0.09315812972761485),
    @("MuSCs", "Dhh", "Hhip", "ECs", 3, 1, 1.583934333333334, 4.751803000000001, 0.4895300409493609, 0.4895300409493609, 3, 1, 0.3309496666666666, 0.9928489999999999, 0.02814339759365984, 0.02814339759365983, 0.5242025396385556, 4.717822856747, 0.01377703857647844, 0.01377703857647844),
    @("MuSCs", "Dhh", "Hhip", "FAPs", 3, 1, 1.583934333333334, 4.751803000000001, 0.4895300409493609, 0.4895300409493609, 3, 1, 6.605175333333333, 19.815526, 0.561692892620634, 0.561692892620634, 10.46216398815311, 94.159475893378, 0.274965544725544, 0.274965544725544),
    @("MuSCs", "Dhh", "Hhip", "MuSCs", 3, 1, 1.583934333333334, 4.751803000000001, 0.4895300409493609, 0.4895300409493609, 3, 1, 4.823282, 14.469846, 0.410163709785706, 0.410163709785706, 7.639761959148667, 68.75785763233802, 0.2007874576473385, 0.2007874576473385)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
